# Adiciona exportação com xlsxwriter e remove indicador sem localização
# Preenche a coluna SEÇÃO (H) com a categoria de cada produto, de acordo
# com a localização/seção da loja a que pertence.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Intervalos de linhas de produto (linha 1 é o cabeçalho) e a seção
# correspondente a cada bloco de produtos.
$ws.Range("H2:H27").Value = "LIMPEZA"
$ws.Range("H28:H77").Value = "MERCEARIA DOCE"
$ws.Range("H78:H143").Value = "MERCEARIA SALGADA"
$ws.Range("H144:H152").Value = "PERFUMARIA"
$ws.Range("H153:H156").Value = "PERFUMARIA INFANTIL"
$ws.Range("H157:H167").Value = "PETSHOP"

# Ajusta a posição da janela/seleção ativa, como deixado após o
# preenchimento manual dos dados (rolagem até a última seção editada).
$ws.Activate()
$ws.Range("H177").Select()
$excel.ActiveWindow.ScrollRow = 157
$excel.ActiveWindow.ScrollColumn = 1
